$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two existing data rows (old rows 21 and 22) so that after
# inserting a new row at the top the sheet still ends at row 21.
$ws.Rows.Item(22).EntireRow.Delete()
$ws.Rows.Item(21).EntireRow.Delete()

# Insert a new row at row 2 (just below the header), pushing the remaining
# data rows down by one.
$ws.Rows.Item(2).EntireRow.Insert()

# Populate the newly inserted row with the new data point.
$ws.Range("A2").Value = 0.0209221355617046
$ws.Range("B2").Value = -0.0198531206697225
$ws.Range("C2").Value = -0.0239764600992202

# The insert operation copies formatting from the row above (the bold
# header); restore plain/default formatting on the new data row to match
# the other data rows.
$ws.Range("A2:C2").ClearFormats()
